{"js": "// Apply hybrid bold + color (#2C3E50) highlighting to quantitative\n// impact metrics (percentages, dollar amounts, large numbers) across\n// the achievements / work-experience bullet paragraphs.\n//\n// Strategy: locate each target paragraph by a unique substring of its\n// text, then within that paragraph's range search (in left-to-right\n// order) for each metric token and bold + color it. Because Word\n// (and this engine) automatically splits a run when a sub-range's\n// formatting is changed, this reproduces the exact run-splitting\n// structure described by the diff without having to hand-build runs.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: a substring that uniquely identifies the target\n// paragraph, followed by the list of metric substrings (in the order\n// they appear in that paragraph) that must become bold + colored.\nconst targets = [\n  {\n    paragraph: \"Discovered systematic race coding errors affecting all Black and Asian-American voters\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    paragraph: \"Utilized advanced sampling methods to decrease survey margin of error\",\n    metrics: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    paragraph: \"Trigonometric algorithm for boundary estimation reduced mapping costs\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    paragraph: \"Built real-time FEC analysis systems using Python, Pandas and PySpark\",\n    metrics: [\"$2\"],\n  },\n  {\n    paragraph: \"Modernized legacy ETL processes by implementing dbt and PySpark workflows\",\n    metrics: [\"57%\"],\n  },\n  {\n    paragraph: \"Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs\",\n    metrics: [\"73.5%\"],\n  },\n  {\n    paragraph: \"$4.7M savings enabled nonprofit access\",\n    metrics: [\"$4.7M\"],\n  },\n  {\n    paragraph: \"Platform impact: Built redistricting system serving\",\n    metrics: [\"12,847\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const target of targets) {\n  // Find the paragraph whose text contains the identifying substring.\n  let paragraph = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(target.paragraph) !== -1) {\n      paragraph = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!paragraph) {\n    throw new Error(\"Paragraph not found for: \" + target.paragraph);\n  }\n\n  // Track how many times each metric string has already been matched\n  // within this paragraph, so repeated tokens (e.g. \"73.5%\" appearing\n  // once per paragraph here, but generically could repeat) are each\n  // bolded independently and only once.\n  const seenCounts = {};\n\n  for (const metric of target.metrics) {\n    const occurrenceIndex = seenCounts[metric] || 0;\n    seenCounts[metric] = occurrenceIndex + 1;\n\n    const found = paragraph.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    if (found.items.length <= occurrenceIndex) {\n      throw new Error(\n        \"Metric '\" + metric + \"' occurrence \" + occurrenceIndex +\n        \" not found in paragraph: \" + target.paragraph\n      );\n    }\n\n    const range = found.items[occurrenceIndex];\n    range.font.bold = true;\n    range.font.color = HIGHLIGHT_COLOR;\n  }\n  await context.sync();\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative\n# impact metrics (percentages, dollar amounts, large numbers) across\n# the achievements / work-experience bullet paragraphs.\n#\n# Strategy: for each target paragraph (identified by a unique\n# substring), duplicate its Range for every metric substring it\n# contains (in left-to-right order), use Find.Execute to collapse\n# that duplicated range onto the metric text, then set Font.Bold /\n# Font.Color on the collapsed range. Word splits the run around the\n# found text automatically, matching the run layout in the diff.\n\nfunction Get-WdColor([string]$hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    # Word COM colors are packed as 0x00BBGGRR (OLE COLORREF).\n    return $r + ($g * 256) + ($b * 65536)\n}\n\n$highlightColor = Get-WdColor \"2C3E50\"\n\n#  Note: single-quoted (literal) strings are used for any metric that\n# contains a '$' so it is never mistaken for PowerShell variable\n# interpolation, regardless of host interpreter.\n$targets = @(\n    @{ Match = 'Discovered systematic race coding errors affecting all Black and Asian-American voters'; Metrics = @('23%', '64%') },\n    @{ Match = 'Utilized advanced sampling methods to decrease survey margin of error'; Metrics = @('\u00b14.2%', '\u00b12.1%', '71%', '87%') },\n    @{ Match = 'Trigonometric algorithm for boundary estimation reduced mapping costs'; Metrics = @('73.5%', '$4.7M') },\n    @{ Match = 'Built real-time FEC analysis systems using Python, Pandas and PySpark'; Metrics = @('$2') },\n    @{ Match = 'Modernized legacy ETL processes by implementing dbt and PySpark workflows'; Metrics = @('57%') },\n    @{ Match = 'Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs'; Metrics = @('73.5%') },\n    @{ Match = '$4.7M savings enabled nonprofit access'; Metrics = @('$4.7M') },\n    @{ Match = 'Platform impact: Built redistricting system serving'; Metrics = @('12,847') }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($target in $targets) {\n    $paragraph = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like (\"*\" + $target.Match + \"*\")) {\n            $paragraph = $p\n            break\n        }\n    }\n\n    if ($null -eq $paragraph) {\n        Write-Output (\"Paragraph not found for: \" + $target.Match)\n        continue\n    }\n\n    foreach ($metric in $target.Metrics) {\n        $searchRange = $paragraph.Range.Duplicate\n        $find = $searchRange.Find\n        $find.Text = $metric\n        $find.MatchCase = $true\n        $find.MatchWildcards = $false\n        $found = $find.Execute()\n\n        if ($found) {\n            $searchRange.Font.Bold = $true\n            $searchRange.Font.Color = $highlightColor\n        } else {\n            Write-Output (\"Metric not found: \" + $metric + \" in \" + $target.Match)\n        }\n    }\n}\n"}
